$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-03-01 03:48:38"
$ws.Range("N2").Value = "-1.2 °C 3:24 TU"
$ws.Range("E3").Value = "2026-03-01 03:48:40"
$ws.Range("E4").Value = "2026-03-01 03:48:43"
$ws.Range("E5").Value = "2026-03-01 03:48:45"
$ws.Range("H5").Value = "'94%"
$ws.Range("N5").Value = "-4.4 °C 3:25 TU"
$ws.Range("O5").Value = "-3.4 °C"
$ws.Range("E6").Value = "2026-03-01 03:48:48"
$ws.Range("H6").Value = "'84%"
$ws.Range("J6").Value = "1025.5 hPa"
$ws.Range("N6").Value = "9.0 °C 3:26 TU"
$ws.Range("E7").Value = "2026-03-01 03:48:51"
$ws.Range("J7").Value = "1025.5 hPa"
$ws.Range("N7").Value = "13.1 °C 3:28 TU"
$ws.Range("E8").Value = "2026-03-01 03:48:53"
$ws.Range("J8").Value = "1025.6 hPa"
$ws.Range("E9").Value = "2026-03-01 03:48:56"
$ws.Range("E10").Value = "2026-03-01 03:48:58"
$ws.Range("O10").Value = "6.9 °C"
$ws.Range("E11").Value = "2026-03-01 03:49:01"
$ws.Range("E12").Value = "2026-03-01 03:49:03"
$ws.Range("E13").Value = "2026-03-01 03:49:06"
$ws.Range("H13").Value = "'90%"
$ws.Range("J13").Value = "1026.3 hPa"
$ws.Range("N13").Value = "4.1 °C 3:29 TU"
$ws.Range("E14").Value = "2026-03-01 03:49:09"
$ws.Range("H14").Value = "'97%"
$ws.Range("L14").Value = "8.3 km/h - 290º 3:28 TU"
$ws.Range("N14").Value = "10.3 °C 3:28 TU"
$ws.Range("O14").Value = "11.2 °C"
$ws.Range("E15").Value = "2026-03-01 03:49:11"
$ws.Range("O15").Value = "8.2 °C"
$ws.Range("E16").Value = "2026-03-01 03:49:13"
$ws.Range("H16").Value = "'85%"
$ws.Range("N16").Value = "-5.3 °C 3:11 TU"
$ws.Range("O16").Value = "-4.6 °C"
$ws.Range("E17").Value = "2026-03-01 03:49:16"
$ws.Range("L17").Value = "7.6 km/h - 258º 3:06 TU"
$ws.Range("N17").Value = "1.1 °C 3:01 TU"
$ws.Range("E18").Value = "2026-03-01 03:49:18"
$ws.Range("J18").Value = "1025.8 hPa"
$ws.Range("L18").Value = "5.0 km/h - 142º 3:04 TU"
$ws.Range("M18").Value = "8.1 °C 3:19 TU"
$ws.Range("O18").Value = "7.0 °C"
$ws.Range("E19").Value = "2026-03-01 03:49:21"
$ws.Range("N19").Value = "5.9 °C 3:27 TU"
$ws.Range("E20").Value = "2026-03-01 03:49:23"
$ws.Range("L20").Value = "10.1 km/h - 228º 3:20 TU"
$ws.Range("N20").Value = "-3.5 °C 3:26 TU"
$ws.Range("O20").Value = "-2.8 °C"
$ws.Range("E21").Value = "2026-03-01 03:49:26"
$ws.Range("N21").Value = "6.3 °C 3:28 TU"
$ws.Range("O21").Value = "6.7 °C"
$ws.Range("E22").Value = "2026-03-01 03:49:28"
$ws.Range("L22").Value = "9.0 km/h - 343º 3:26 TU"
$ws.Range("N22").Value = "-5.9 °C 3:23 TU"
$ws.Range("O22").Value = "-5.0 °C"
$ws.Range("E23").Value = "2026-03-01 03:49:31"
$ws.Range("L23").Value = "20.9 km/h - 280º 3:04 TU"
$ws.Range("E24").Value = "2026-03-01 03:49:34"
$ws.Range("O24").Value = "4.2 °C"
$ws.Range("E25").Value = "2026-03-01 03:49:37"
$ws.Range("L25").Value = "16.2 km/h - 249º 3:27 TU"
$ws.Range("N25").Value = "-2.5 °C 3:29 TU"
$ws.Range("O25").Value = "-1.9 °C"
$ws.Range("E26").Value = "2026-03-01 03:49:39"
$ws.Range("H26").Value = "'97%"
$ws.Range("J26").Value = "1025.9 hPa"
$ws.Range("N26").Value = "2.5 °C 3:29 TU"
$ws.Range("E27").Value = "2026-03-01 03:49:42"
$ws.Range("O27").Value = "-1.3 °C"
$ws.Range("E28").Value = "2026-03-01 03:49:44"
$ws.Range("E29").Value = "2026-03-01 03:49:47"
$ws.Range("E30").Value = "2026-03-01 03:49:50"
$ws.Range("H30").Value = "'81%"
$ws.Range("J30").Value = "1025.5 hPa"
$ws.Range("O30").Value = "10.0 °C"
$ws.Range("E31").Value = "2026-03-01 03:49:52"
$ws.Range("N31").Value = "10.7 °C 3:20 TU"
$ws.Range("E32").Value = "2026-03-01 03:49:55"
$ws.Range("O32").Value = "1.7 °C"
$ws.Range("E33").Value = "2026-03-01 03:49:58"
$ws.Range("H33").Value = "'91%"
$ws.Range("J33").Value = "1025.7 hPa"
$ws.Range("E34").Value = "2026-03-01 03:50:00"
$ws.Range("N34").Value = "-0.3 °C 3:10 TU"
$ws.Range("O34").Value = "-0.1 °C"
$ws.Range("E35").Value = "2026-03-01 03:50:03"
$ws.Range("E36").Value = "2026-03-01 03:50:06"
$ws.Range("H36").Value = "'79%"
$ws.Range("L36").Value = "9.7 km/h - 154º 3:23 TU"
$ws.Range("O36").Value = "9.3 °C"
$ws.Range("E37").Value = "2026-03-01 03:50:08"
$ws.Range("H37").Value = "'97%"
$ws.Range("J37").Value = "1026.7 hPa"
$ws.Range("L37").Value = "5.8 km/h - 64º 3:05 TU"
$ws.Range("N37").Value = "6.1 °C 3:11 TU"
$ws.Range("E38").Value = "2026-03-01 03:50:11"
$ws.Range("M38").Value = "9.5 °C 3:03 TU"
$ws.Range("O38").Value = "8.7 °C"
$ws.Range("E39").Value = "2026-03-01 03:50:13"
$ws.Range("E40").Value = "2026-03-01 03:50:16"
$ws.Range("J40").Value = "1025.6 hPa"
$ws.Range("N40").Value = "6.8 °C 3:05 TU"
$ws.Range("O40").Value = "7.4 °C"
$ws.Range("E41").Value = "2026-03-01 03:50:19"
$ws.Range("H41").Value = "'86%"
$ws.Range("J41").Value = "1025.5 hPa"
$ws.Range("M41").Value = "12.0 °C 3:12 TU"
$ws.Range("E42").Value = "2026-03-01 03:50:21"
$ws.Range("H42").Value = "'81%"
$ws.Range("N42").Value = "7.3 °C 3:22 TU"
$ws.Range("O42").Value = "9.5 °C"
$ws.Range("E43").Value = "2026-03-01 03:50:24"
$ws.Range("N43").Value = "8.4 °C 3:00 TU"
$ws.Range("O43").Value = "8.7 °C"
$ws.Range("E44").Value = "2026-03-01 03:50:27"
$ws.Range("N44").Value = "-3.0 °C 3:15 TU"
$ws.Range("O44").Value = "-2.4 °C"
$ws.Range("E45").Value = "2026-03-01 03:50:29"
$ws.Range("J45").Value = "1027.1 hPa"
$ws.Range("N45").Value = "3.5 °C 3:12 TU"
$ws.Range("E46").Value = "2026-03-01 03:50:32"
$ws.Range("O46").Value = "7.5 °C"
